$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.510.34'
$ws.Range("E2").Value = '  +1.61%  '
$ws.Range("D3").Value = '1.573.13'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -1.73%  '
$ws.Range("D5").Value = '211.47'
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("E7").Value = '  -1.63%  '
$ws.Range("D8").Value = '22.87'
$ws.Range("E8").Value = '  +3.30%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").Value = '0.0595'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").Value = '1.797.95'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '1.554.74'
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '27.501.89'
$ws.Range("D17").Value = '62.40'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '226.32'
$ws.Range("E18").Value = '  +4.73%  '
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = '0.0₃0705'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("E21").Value = '  -1.68%  '
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").Value = '9.41'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").Value = '150.34'
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D26").Value = '15.17'
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("D29").Value = '0.993'
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("E30").Value = '  +0.69%  '
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").Value = '1.457.09'
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("E34").Value = '  -2.11%  '
$ws.Range("E35").Value = '  +3.48%  '
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("D41").Value = '2.36'
$ws.Range("E41").Value = '  -3.46%  '
$ws.Range("E42").Value = '  -3.17%  '
$ws.Range("D43").Value = '0.992'
$ws.Range("E43").Value = '  -1.84%  '
$ws.Range("E44").Value = '  +6.83%  '
$ws.Range("D45").Value = '0.974'
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("D46").Value = '64.14'
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("D47").Value = '1.709.82'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").Value = '87.01'
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("E49").Value = '  +2.26%  '
$ws.Range("D50").Value = '0.0528'
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").Value = '0.0947'
$ws.Range("E51").Value = '  -1.85%  '
